$wb = $excel.ActiveWorkbook

# Rename "measurement" -> "observation"
$wsObs = $wb.Worksheets.Item("measurement")
$wsObs.Name = "observation"

# Add a new worksheet "variance" right after "observation"
$wsVar = $wb.Worksheets.Add($null, $wsObs)
$wsVar.Name = "variance"

# Move column B ("Variance") values out of "observation" into column A of "variance"
$lastRow = $wsObs.Cells.Item($wsObs.Rows.Count, 2).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $val = $wsObs.Cells.Item($r, 2).Value2
    $wsVar.Cells.Item($r, 1).Value = $val
}

# Remove column B from "observation" now that its data lives in "variance"
$wsObs.Columns.Item(2).ClearContents()

# Column width for new sheet's column A (matches old column B width in observation,
# closest reachable value given the engine's pixel-quantized ColumnWidth conversion)
$wsVar.Columns.Item(1).ColumnWidth = 7.5

# Restore selections similar to the target state
[void]$wsObs.Activate()
[void]$wsObs.Range("D11").Select()

[void]$wsVar.Activate()
[void]$wsVar.Range("D15").Select()

[void]$wsObs.Activate()
